# Commit: "Natmi following Dr Hou advice"
# Adds a new "ECs" sending/target cluster to the Wnt1 -> Ryk ligand-receptor
# results sheet. The table grows from 2 rows (FAPs only) to 6 rows, covering
# every combination of {ECs, FAPs} as sending cluster and {ECs, FAPs, sCs} as
# target cluster, with refreshed values in every derived/statistical column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Row 2..7 values, column order matches $cols (A:T)
$newRows = @(
  @("ECs", "Wnt1", "Ryk", "ECs", 1, 0.3333333333333333, 0.027123, 0.081369, 0.07131444737854614, 0.07131444737854616, 2, 0.6666666666666666, 6.245403666666667, 18.736211, 0.1461562881451252, 0.1461562881451252, 0.169394083651, 1.524546752859, 0.01042305491996916, 0.01042305491996916),
  @("ECs", "Wnt1", "Ryk", "FAPs", 1, 0.3333333333333333, 0.027123, 0.081369, 0.07131444737854614, 0.07131444737854616, 3, 1, 21.552384, 64.657152, 0.5043735544158399, 0.5043735544158399, 0.584565311232, 5.261087801087999, 0.03596912130551869, 0.0359691213055187),
  @("ECs", "Wnt1", "Ryk", "sCs", 1, 0.3333333333333333, 0.027123, 0.081369, 0.07131444737854614, 0.07131444737854616, 3, 1, 14.93320766666667, 44.799623, 0.3494701574390349, 0.3494701574390349, 0.4050333915429999, 3.645300523887, 0.02492227115305829, 0.0249222711530583),
  @("FAPs", "Wnt1", "Ryk", "ECs", 3, 1, 0.3532066666666667, 1.05962, 0.9286855526214538, 0.9286855526214538, 2, 0.6666666666666666, 6.245403666666667, 18.736211, 0.1461562881451252, 0.1461562881451252, 2.205918211091111, 19.85326389982, 0.135733233225156, 0.135733233225156),
  @("FAPs", "Wnt1", "Ryk", "FAPs", 3, 1, 0.3532066666666667, 1.05962, 0.9286855526214538, 0.9286855526214538, 3, 1, 21.552384, 64.657152, 0.5043735544158399, 0.5043735544158399, 7.61244571136, 68.51201140223999, 0.4684044331103212, 0.4684044331103212),
  @("FAPs", "Wnt1", "Ryk", "sCs", 3, 1, 0.3532066666666667, 1.05962, 0.9286855526214538, 0.9286855526214538, 3, 1, 14.93320766666667, 44.799623, 0.3494701574390349, 0.3494701574390349, 5.274508502584444, 47.47057652326, 0.3245478862859766, 0.3245478862859766)
)

$startRow = 2
for ($r = 0; $r -lt $newRows.Length; $r++) {
  $rowVals = $newRows[$r]
  $rowNum = $startRow + $r
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $addr = $cols[$i] + $rowNum
    $ws.Range($addr).Value = $rowVals[$i]
  }
}
